# "change excel2json to many" - refresh the monitor report figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 - 總覽 (overview): the values in this sheet are stored as text
# labels (e.g. "31.65", "-3.64%") rather than numbers, so force the
# cells to Text format before writing to stop Excel from re-parsing the
# numeric-looking strings back into numbers/percentages.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C5:D6").NumberFormat = "@"
$ws1.Range("C5").Value = "4861.4$"
$ws1.Range("D5").Value = "-1.20%"
$ws1.Range("C6").Value = "31.66"
$ws1.Range("D6").Value = "+0.27%"

$ws1.Range("C11").NumberFormat = "@"
$ws1.Range("C11").Value = "158.39%"
$ws1.Range("E11").NumberFormat = "@"
$ws1.Range("E11").Value = "140.85%"

# ---------------------------------------------------------------------
# Sheet 2 - 詳細數據 (detail stats): P/C ratio 5-day average, also text.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B21").NumberFormat = "@"
$ws2.Range("B21").Value = "140.85%"

# ---------------------------------------------------------------------
# Sheet 3 - 個股籌碼 (per-stock chips): drop the broker-level columns
# (P..U) for every data row and mark the data-source column as N/A;
# also refresh a handful of price/volume figures.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("P4:U19").Clear()
$ws3.Range("V4:V19").Value = "N/A"

# Row 15 (3081 聯亞)
$ws3.Range("C15").Value = 1000
$ws3.Range("D15").Value = -2.44
$ws3.Range("E15").Value = 905

# Row 16 (3260 威剛) - % change flips from positive to negative, so the
# conditional red/green font swaps too.
$ws3.Range("C16").Value = 306.5
$ws3.Range("D16").Value = -8.1
$ws3.Range("D16").Font.Color = 32768
$ws3.Range("E16").Value = 26486

# Row 17 (3265 台新科)
$ws3.Range("C17").Value = 132.5
$ws3.Range("D17").Value = -2.57
$ws3.Range("E17").Value = 1591

# Row 18 (4979 華星光)
$ws3.Range("C18").Value = 318.5
$ws3.Range("D18").Value = -0.16
$ws3.Range("E18").Value = 20497
